$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 60 - 四方坪站 (shared string index 2)
$ws.Range("A60").Value = 45930
$ws.Range("B60").Value = "四方坪站"
$ws.Range("C60").Formula = "=19745/127"
$ws.Range("D60").Formula = "=C60/(24*60)"
$ws.Range("E60").Formula = "=10978.97/127"
$ws.Range("F60").Formula = "=3797.81/127"
$ws.Range("G60").Formula = "=10978.91/(19745/60)"
$ws.Range("H60").Formula = "=475/127"

# Row 61 - 高岭站 (shared string index 3)
$ws.Range("A61").Value = 45930
$ws.Range("B61").Value = "高岭站"
$ws.Range("C61").Formula = "=7675/36"
$ws.Range("D61").Formula = "=C61/(24*60)"
$ws.Range("E61").Formula = "=5463.25/36"
$ws.Range("F61").Formula = "=1342.68/36"
$ws.Range("G61").Formula = "=5463.25/(7675/60)"
$ws.Range("H61").Formula = "=218/36"

$ws.Range("J61").Select()
